$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Insert new row 5: complete / COMPLETE  (pushes old rows 5.. down by one)
# ---------------------------------------------------------------------------
$ws.Rows.Item(5).Insert()
$ws.Range("A5").Value = "complete"
$ws.Range("B5").Value = "COMPLETE"
$ws.Range("B5").WrapText = $true

# ---------------------------------------------------------------------------
# 2) Insert new row 33: total_score / Total Score:
#    (row 32 is now "total" / "Total:" after the first insert)
# ---------------------------------------------------------------------------
$ws.Rows.Item(33).Insert()
$ws.Range("A33").Value = "total_score"
$ws.Range("B33").Value = "Total Score:"

# ---------------------------------------------------------------------------
# 3) Add VoiceDuration (column C) to two existing rows
# ---------------------------------------------------------------------------
$ws.Range("C52").Value = 2.5
$ws.Range("C53").Value = 2.5

# ---------------------------------------------------------------------------
# 4) Merge/replace the two level_0_intro_1_4 / level_0_intro_1_5 rows into one
# ---------------------------------------------------------------------------
$ws.Range("B62").Value = "A unit cube's measurement can also be changed based on specific needs. In our case, one unit cube equals to one cubic feet."
$ws.Rows.Item(63).Delete()

# ---------------------------------------------------------------------------
# 5) Append new rows at the end (after the former row 63, now row 64)
# ---------------------------------------------------------------------------
$newRows = @(
    @("level_0_end_1", "If you count the number of unit cubes placed on the ground, it tells you the volume of the object."),
    @("level_0_end_2", "In this case, this object is made up of 16 unit cubes, where each cube's volume is 1 cubic feet. Therefore the volume of the object is 16 cubic feet."),
    @("level_0_end_3", "Anyhow, it's time to build!"),
    @("level_1_intro_0_1", "For this objective, we will have to stack more than one layer of cubes."),
    @("level_1_intro_0_2", "In order to increase the stack while expanding, simply highlight the top surface and drag upwards."),
    @("level_2_intro_0_1", "On this level, the unit cube's sides are half a foot with a volume of one-eighth cubic feet."),
    @("level_2_intro_0_2", "Since we are dealing with measuring the volume using unit cubes, first you need to compute the number of unit cubes."),
    @("level_2_intro_0_3", "Then simply multiply the number of unit cubes with the unit cube's volume to get the correct result."),
    @("level_2_intro_0_4", "So for this level, we will need 8 unit cubes. Multiply 8 with one-eighth to convert it to the correct volume, which is 1 cubic foot."),
    @("level_2_unit_formula", "2 Units x 2 Units x 2 Units = 8 Cubic Units"),
    @("level_2_unit_volume", "Cubic Unit Volume = 1/8 ft³"),
    @("level_2_unit_to_volume", "8 Cubic Units x 1/8 ft³ = 1 ft³"),
    @("level_4_intro_0_1", "On this level, you will need to place two groups of unit cubes to get the required amount of volume."),
    @("level_4_intro_0_2", "Volumes can be added together in any shape or form, so long as they are in the same measurement."),
    @("level_4_intro_0_3", "In order for volumes to be added together make sure they are placed adjacently on any of the sides."),
    @("level_6_intro_0_1", "On this level, one of the objectives have a restricted height."),
    @("level_6_intro_0_2", "You won't be able to expand this particular material's height beyond the limit."),
    @("end_title", "CONGRATULATIONS", 1.5),
    @("end_detail", "You have successfully given all the frogs a loving home!", 3),
    @("end_detail_2", "Thank you for playing!", 2)
)

$r = 65
foreach ($row in $newRows) {
    $ws.Range("A$r").Value = $row[0]
    $ws.Range("B$r").Value = $row[1]
    if ($row.Length -gt 2) {
        $ws.Range("C$r").Value = $row[2]
    }
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# View tweaks
# ---------------------------------------------------------------------------
$ws.Range("A33").Select()
$excel.ActiveWindow.ScrollRow = 22
